# summer 24 week 5 inputs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.3
$ws.Range("B3").Value = 1.54
$ws.Range("E3").Value = 1.33
$ws.Range("C4").Value = 1.43
$ws.Range("E4").Value = 1.22
$ws.Range("F4").Value = 1.11
$ws.Range("C5").Value = 1.33
$ws.Range("F5").Value = 1.05
$ws.Range("D6").Value = 1.5
$ws.Range("E6").Value = 1.32
